$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

function Insert-ItalicParagraphAfter($anchorText, $newText) {
    $rng = $d.Content
    $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $rng.InsertParagraphAfter()
    $newRng = $d.Range($rng.End + 1, $rng.End + 1)
    $newRng.InsertAfter($newText)
    $newRng.Font.Italic = $true
}

# 1. Heading3 paragraph gets the English course title.
$rng = $d.Content
$rng.Find.Execute("LOQ4241 -  Sistemas de Apoio à Decisão", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$heading3 = $d.Paragraphs(2).Range
$heading3.InsertAfter("Decision  Support Systems")

# 2. Update activation date.
Replace-Text "Ativação: 01/01/2016" "Ativação: 01/01/2024"

# 3. Insert English objective paragraph (italic) after the Portuguese one.
Insert-ItalicParagraphAfter "Apresentar conceitos, ferramentas e métodos para o auxílio à tomada de decisão." "Provide theory, tools and methods for supporting decision-making."

# 4. Update the responsible professor.
Replace-Text "5840917 - Fabricio Maciel Gomes" "3295113 - José Eduardo Holler Branco"

# 5. Update the summarized program text, then add its English translation.
Replace-Text "Teoria da Decisão; Estruturação, Decisão sem Risco e sem Incerteza; Decisão com Múltiplos Cenários ou Múltiplos Critérios; Decisão com Incerteza; Sistemas de Auxílio à Decisão e Sistemas Especialistas." "Teoria da Decisão; Planejamento de sistemas de apoio à decisão; Decisão com Múltiplos Cenários ou Múltiplos Critérios; Decisão com Incerteza; Sistemas de Auxílio à Decisão."
Insert-ItalicParagraphAfter "Teoria da Decisão; Planejamento de sistemas de apoio à decisão; Decisão com Múltiplos Cenários ou Múltiplos Critérios; Decisão com Incerteza; Sistemas de Auxílio à Decisão." "Decision Theory; Decision support systems planning; Decision with Multiple Scenarios or Multiple Criteria; Decision with Uncertainty; Decision Support Systems."

# 6. Replace the detailed program list (was multiple runs joined by line breaks) with one run, then add its English translation.
$rng = $d.Content
$rng.Find.Execute("1.Teoria da Decisão", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$progPara = $rng.Paragraphs(1)
$progPara.Range.Text = "i) Teoria da Decisão; ii) Estruturação de modelos de decisão; iii) Decisão com múltiplos cenários ou múltiplos critérios; iv) Decisão com incerteza; e v) Sistemas de auxílio à decisão e sistemas especialistas."
Insert-ItalicParagraphAfter "i) Teoria da Decisão; ii) Estruturação de modelos de decisão; iii) Decisão com múltiplos cenários ou múltiplos critérios; iv) Decisão com incerteza; e v) Sistemas de auxílio à decisão e sistemas especialistas." "i) Decision Theory; ii) Decision support systems planning; iii) Decision with Multiple Scenarios or Multiple Criteria; iv) Decision with Uncertainty; and v) Decision Support Systems."

# 7. Evaluation method.
Replace-Text "Aulas expositivas teóricas, aulas práticas, aulas de exercícios." "Provas, trabalhos em grupo, exercícios individuais e seminários."

# 8. Evaluation criterion.
Replace-Text "A Nota Final do aluno será determinada segundo a seguinte equação: Nota Final = (Prova- Bimestral-1*0,4) + (Prova-Bimestral-2*0,4) + (Trabalho*0,2)" "Média das atividades avaliativas."

# 9. Recovery rule.
Replace-Text "Média aritmética da nota final obtida pelo aluno durante o semestre e da nota obtida na Prova de Recuperação." "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."

# 10. Replace the bibliography (was multiple runs joined by line breaks) with one combined run.
$rng = $d.Content
$rng.Find.Execute("1.ENSSLIN, L.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bibPara = $rng.Paragraphs(1)
$bibPara.Range.Text = "FURTADO, N.; KAWAMOTO, E. Avaliação de Projetos de Transporte. São Carlos: Serviço Gráfico EESC-USP, 2002. 254 p.POWER, D. J. Decision Support Systems. London: Quorum Books, 2002. 251 p.GOMES, L. F. A. M.; GOMES, C. F. S.; ALMEIDA, A. T, Tomada de Decisão Gerencial: enfoque multicritério, São Paulo: Atlas, 2002.SHIMIZU, T., Decisão nas Organizações: introdução aos problemas de decisão encontrados nas organizações e nos sistemas de apoio à decisão, São Paulo: Atlas, 2001.DEVLIN, G. (ed.). Decision Support Systems: advances in. Zagreb: Intech, 2010. 342 p.GARCÍA-DÍAZ, V. Algorithms in Decision Support Systems. Basel: MDPI, 2020. 147 p."

Write-Host "All edits applied"
